# Auto-generated edit script for 北京-漫展信息.xlsx
# Applies the scraped-data refresh described in the commit diff:
#  - 展览 (sheet1): "想去人数" (F) count refresh for many rows
#  - 演出 (sheet2): row 4 (2024-08-19 event) removed; rows 5-16 shift up to 4-15
#  - 本地生活 (sheet3): F4 count refresh
#  - 全部类型 (sheet4): "想去人数" (F) count refresh for many rows

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- 展览: refresh "想去人数" (column F) ----
$ws1.Cells.Item(4, 6).Value = 148
$ws1.Cells.Item(5, 6).Value = 61
$ws1.Cells.Item(6, 6).Value = 3759
$ws1.Cells.Item(7, 6).Value = 221
$ws1.Cells.Item(8, 6).Value = 2519
$ws1.Cells.Item(9, 6).Value = 65
$ws1.Cells.Item(10, 6).Value = 3022
$ws1.Cells.Item(11, 6).Value = 1394
$ws1.Cells.Item(12, 6).Value = 528
$ws1.Cells.Item(13, 6).Value = 2276
$ws1.Cells.Item(15, 6).Value = 113
$ws1.Cells.Item(16, 6).Value = 82
$ws1.Cells.Item(17, 6).Value = 433
$ws1.Cells.Item(19, 6).Value = 190
$ws1.Cells.Item(21, 6).Value = 297
$ws1.Cells.Item(22, 6).Value = 328
$ws1.Cells.Item(23, 6).Value = 641
$ws1.Cells.Item(24, 6).Value = 1381
$ws1.Cells.Item(25, 6).Value = 36
$ws1.Cells.Item(26, 6).Value = 1288
$ws1.Cells.Item(29, 6).Value = 240
$ws1.Cells.Item(30, 6).Value = 19
$ws1.Cells.Item(31, 6).Value = 4164
$ws1.Cells.Item(32, 6).Value = 3792
$ws1.Cells.Item(33, 6).Value = 65
$ws1.Cells.Item(35, 6).Value = 1103
$ws1.Cells.Item(36, 6).Value = 450
$ws1.Cells.Item(37, 6).Value = 8
$ws1.Cells.Item(38, 6).Value = 1304
$ws1.Cells.Item(39, 6).Value = 144
$ws1.Cells.Item(40, 6).Value = 123
$ws1.Cells.Item(41, 6).Value = 87

# ---- 本地生活: refresh "想去人数" (column F) ----
$ws3.Cells.Item(4, 6).Value = 2234

# ---- 全部类型: refresh "想去人数" (column F) ----
$ws4.Cells.Item(7, 6).Value = 148
$ws4.Cells.Item(8, 6).Value = 61
$ws4.Cells.Item(10, 6).Value = 3759
$ws4.Cells.Item(11, 6).Value = 221
$ws4.Cells.Item(12, 6).Value = 2519
$ws4.Cells.Item(13, 6).Value = 65
$ws4.Cells.Item(14, 6).Value = 3022
$ws4.Cells.Item(15, 6).Value = 528
$ws4.Cells.Item(16, 6).Value = 2276
$ws4.Cells.Item(18, 6).Value = 113
$ws4.Cells.Item(19, 6).Value = 82
$ws4.Cells.Item(20, 6).Value = 433
$ws4.Cells.Item(22, 6).Value = 190
$ws4.Cells.Item(24, 6).Value = 328
$ws4.Cells.Item(25, 6).Value = 641
$ws4.Cells.Item(26, 6).Value = 1381
$ws4.Cells.Item(27, 6).Value = 36
$ws4.Cells.Item(28, 6).Value = 1288
$ws4.Cells.Item(32, 6).Value = 20
$ws4.Cells.Item(33, 6).Value = 4164
$ws4.Cells.Item(34, 6).Value = 3793
$ws4.Cells.Item(35, 6).Value = 65
$ws4.Cells.Item(38, 6).Value = 450
$ws4.Cells.Item(40, 6).Value = 8
$ws4.Cells.Item(43, 6).Value = 1304
$ws4.Cells.Item(44, 6).Value = 144
$ws4.Cells.Item(45, 6).Value = 87
$ws4.Cells.Item(49, 6).Value = 198

# ---- 演出: the 2024-08-19 event (old row 4) was removed upstream.
# Shift rows 5-16 (columns B-I) up into rows 4-15 (column A, the static
# sequence index, is left untouched), then delete the now-duplicated
# last row (16) to shrink the used range back down to A1:I15.
# row 4 (was row 5)
$ws2.Cells.Item(4, 2).Value = '2024-08-23'
$ws2.Cells.Item(4, 3).Value = '北京·《山丘》音乐教父 经典情歌金曲翻唱演唱会'
$ws2.Cells.Item(4, 4).Value = '大江胡同121号2幢负1层 北京门空间 TheDoorLiveHouse'
$ws2.Cells.Item(4, 5).Value = '2024.08.23 19:30-08.23 21:00'
$ws2.Cells.Item(4, 6).Value = 1
$ws2.Cells.Item(4, 7).Value = 98
$ws2.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89358'
$ws2.Cells.Item(4, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/noqwx8Qu1721116074567.jpeg'
# row 5 (was row 6)
$ws2.Cells.Item(5, 2).Value = '2024-08-24'
$ws2.Cells.Item(5, 3).Value = '北京·最后的莫西干人——亚历桑德罗&丛林回响乐队印第安音乐品鉴会'
$ws2.Cells.Item(5, 4).Value = '亮马桥路40号(近好运街) 北京世纪剧院'
$ws2.Cells.Item(5, 5).Value = '2024.08.24 19:30-08.24 21:00'
$ws2.Cells.Item(5, 6).Value = 18
$ws2.Cells.Item(5, 7).Value = 153
$ws2.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86217'
$ws2.Cells.Item(5, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/BDyblKrJ1716427731729.jpeg'
# row 6 (was row 7)
$ws2.Cells.Item(6, 2).Value = '2024-09-22'
$ws2.Cells.Item(6, 3).Value = '北京·《喜剧奇妙夜》一年一度喜剧大赛编剧团队编创/切西娅剧组演绎'
$ws2.Cells.Item(6, 4).Value = '复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)'
$ws2.Cells.Item(6, 5).Value = '2024.09.22 19:30-09.22 21:00'
$ws2.Cells.Item(6, 6).Value = 3
$ws2.Cells.Item(6, 7).Value = 80
$ws2.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90700'
$ws2.Cells.Item(6, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/tCGETFGW1723613408321.jpeg'
# row 7 (was row 8)
$ws2.Cells.Item(7, 2).Value = '2024-09-22'
$ws2.Cells.Item(7, 3).Value = '北京·次元音浪Million Live⏤番音集结'
$ws2.Cells.Item(7, 4).Value = '学清路38号金码大厦B座 北京想象空间'
$ws2.Cells.Item(7, 5).Value = '2024.09.22 13:00-09.22 16:00'
$ws2.Cells.Item(7, 6).Value = 20
$ws2.Cells.Item(7, 7).Value = 88
$ws2.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90657'
$ws2.Cells.Item(7, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/Fn9CSOmf1723477511986.jpeg'
# row 8 (was row 9)
$ws2.Cells.Item(8, 2).Value = '2024-09-30'
$ws2.Cells.Item(8, 3).Value = '北京·VGL 中国巡演 2024 VIDEO GAME LIVE 魔兽世界音乐会'
$ws2.Cells.Item(8, 4).Value = '西直门外大街135号  北展剧场'
$ws2.Cells.Item(8, 5).Value = '2024.09.30 19:30-09.30 21:30'
$ws2.Cells.Item(8, 6).Value = 14
$ws2.Cells.Item(8, 7).Value = 180
$ws2.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89144'
$ws2.Cells.Item(8, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/5YIwe8lU1720605586333.jpeg'
# row 9 (was row 10)
$ws2.Cells.Item(9, 2).Value = '2024-10-01'
$ws2.Cells.Item(9, 3).Value = '北京·VGL 中国巡演 2024 VIDEO GAMES LIVE 暴雪游戏音乐会'
$ws2.Cells.Item(9, 4).Value = '西直门外大街135号  北展剧场'
$ws2.Cells.Item(9, 5).Value = '2024.10.01 19:30-10.01 21:30'
$ws2.Cells.Item(9, 6).Value = 17
$ws2.Cells.Item(9, 7).Value = '不可售'
$ws2.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89083'
$ws2.Cells.Item(9, 9).Value = '//i1.hdslb.com/bfs/openplatform/202407/yMoDGuXs1720607500874.jpeg'
# row 10 (was row 11)
$ws2.Cells.Item(10, 2).Value = '2024-10-10'
$ws2.Cells.Item(10, 3).Value = '北京·黑白键上的音乐地图——孩子们的钢琴协奏曲之夜'
$ws2.Cells.Item(10, 4).Value = '北新华街1号 北京音乐厅'
$ws2.Cells.Item(10, 5).Value = '2024.10.10 19:30-10.10 21:00'
$ws2.Cells.Item(10, 6).Value = 1
$ws2.Cells.Item(10, 7).Value = 153
$ws2.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86881'
$ws2.Cells.Item(10, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/K3oihoH91717474488019.jpeg'
# row 11 (was row 12)
$ws2.Cells.Item(11, 2).Value = '2024-10-11'
$ws2.Cells.Item(11, 3).Value = '北京·官方唯一授权-周杰伦作品国风音乐会'
$ws2.Cells.Item(11, 4).Value = '西直门外大街135号  北展剧场'
$ws2.Cells.Item(11, 5).Value = '2024.10.11 19:30-10.11 21:00'
$ws2.Cells.Item(11, 6).Value = 14
$ws2.Cells.Item(11, 7).Value = 126
$ws2.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88666'
$ws2.Cells.Item(11, 9).Value = '//i1.hdslb.com/bfs/openplatform/202407/2KgWinEn1720077808243.jpeg'
# row 12 (was row 13)
$ws2.Cells.Item(12, 2).Value = '2024-10-25'
$ws2.Cells.Item(12, 3).Value = '北京·伦敦西区音乐剧明星演唱会-经典版'
$ws2.Cells.Item(12, 4).Value = '西直门外大街135号（北京展览馆内） 北京展览馆剧场'
$ws2.Cells.Item(12, 5).Value = '2024.10.25 19:30-10.26 21:30'
$ws2.Cells.Item(12, 6).Value = 5
$ws2.Cells.Item(12, 7).Value = 144
$ws2.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89359'
$ws2.Cells.Item(12, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/PzPiEKUI1721114840552.jpeg'
# row 13 (was row 14)
$ws2.Cells.Item(13, 2).Value = '2024-10-26'
$ws2.Cells.Item(13, 3).Value = '北京·伦敦西区音乐剧明星演唱会（摇滚版）'
$ws2.Cells.Item(13, 4).Value = '西直门外大街135号（北京展览馆内） 北京展览馆剧场'
$ws2.Cells.Item(13, 5).Value = '2024.10.26 14:30-10.26 16:30'
$ws2.Cells.Item(13, 6).Value = 8
$ws2.Cells.Item(13, 7).Value = 144
$ws2.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89400'
$ws2.Cells.Item(13, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/TYPRpfu21721116217467.jpeg'
# row 14 (was row 15)
$ws2.Cells.Item(14, 2).Value = '2024-10-26'
$ws2.Cells.Item(14, 3).Value = '北京·变形金刚音乐会40周年特别版'
$ws2.Cells.Item(14, 4).Value = '中关村南大街33号国家图书馆北门 国图艺术中心音乐厅'
$ws2.Cells.Item(14, 5).Value = '2024.10.26 19:30-10.26 21:30'
$ws2.Cells.Item(14, 6).Value = 42
$ws2.Cells.Item(14, 7).Value = 266
$ws2.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89970'
$ws2.Cells.Item(14, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/TwvRQI041722150343639.jpeg'
# row 15 (was row 16)
$ws2.Cells.Item(15, 2).Value = '2024-11-30'
$ws2.Cells.Item(15, 3).Value = '北京·花たん 2024 LIVE in Beijing'
$ws2.Cells.Item(15, 4).Value = '复兴路69号院2号136、G23室 Mao Livehouse北京五棵松店'
$ws2.Cells.Item(15, 5).Value = '2024.11.30 14:00-11.30 15:30'
$ws2.Cells.Item(15, 6).Value = 198
$ws2.Cells.Item(15, 7).Value = 380
$ws2.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90341'
$ws2.Cells.Item(15, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/wfGEn3sY1722910561352.jpeg'

$ws2.Rows(16).Delete()

